$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRows = @(5, 6, 8, 10, 11, 12, 15, 19, 20, 21, 22, 23, 28, 32, 33, 34, 36, 37, 38, 39, 41, 42, 45, 46, 47, 49, 50)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "66.982.71"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "3.259.17"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "579.38"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").Value = "176.95"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "3.259.15"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("D11").Value = "6.75"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("D13").Value = "3.831.54"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "28.13"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "67.012.68"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "3.265.58"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "5.84"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").Value = "372.12"
$ws.Range("E21").Value = "  +5.18%  "
$ws.Range("D22").Value = "7.62"
$ws.Range("E22").Value = "  +5.66%  "
$ws.Range("D23").Value = "71.50"
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "3.410.07"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").Value = "5.62"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "22.63"
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").Value = "6.83"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("D37").Value = "167.09"
$ws.Range("E37").Value = "  +7.38%  "
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  +4.22%  "
$ws.Range("D39").Value = "0.856"
$ws.Range("E39").Value = "  +5.44%  "
$ws.Range("E40").Value = "  +9.20%  "
$ws.Range("D41").Value = "27.30"
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("D44").Value = "2.744.92"
$ws.Range("E44").Value = "  +4.80%  "
$ws.Range("D45").Value = "4.36"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "348.61"
$ws.Range("E46").Value = "  +4.79%  "
$ws.Range("D47").Value = "25.07"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").Value = "0.0676"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "0.0280"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("E51").Value = "  +1.50%  "
